$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Order")

# Fill in the "Verkocht" (sold) quantities in column E for rows 2-14
$soldValues = @{
    2  = 2
    3  = 5
    4  = 1
    5  = 7
    6  = 8
    7  = 3
    8  = 4
    9  = 8
    10 = 2
    11 = 1
    12 = 10
    13 = 1
    14 = 0
}

foreach ($row in $soldValues.Keys) {
    $ws.Cells.Item($row, 5).Value = $soldValues[$row]
}

# Add a new row (15) for "big hoops"
$ws.Range("A15").Value = "big hoops"
$ws.Range("B15").Value = 12
$ws.Range("C15").Value = 2.5
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 3
$ws.Range("F15").Formula = "=(B15-C15)*E15"

# Copy formatting from row 14 to row 15 so the new row matches existing styling
$ws.Range("A14:F14").Copy()
$ws.Range("A15:F15").PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Match row height metadata of the other data rows
$ws.Rows.Item(15).RowHeight = 15

# Update selection to match the post-edit state
$ws.Range("E16").Select() | Out-Null
